{"js": "// The \"None\" paragraph (answer to \"Prerequisite tasks\") was incorrectly\n// styled as Heading 1 with an explicit 22-half-point size override. Fix it\n// so it is a normal body-text paragraph (no heading style, no size\n// override), and drop the bookmark that had been placed on it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"None\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  // Demote the paragraph from Heading 1 to Normal (body text).\n  target.style = \"Normal\";\n\n  // Remove the explicit 11pt (half-points value 22) font-size override on\n  // the paragraph's text run(s).\n  const range = target.getRange();\n  range.font.size = null;\n\n  await context.sync();\n}\n\n// The paragraph used to carry a bookmark (\"_aqdfz55armzg\") marking it as\n// the prerequisite-tasks answer; that bookmark is no longer needed once the\n// paragraph is plain body text.\ncontext.document.deleteBookmark(\"_aqdfz55armzg\");\nawait context.sync();\n", "ps1": "# The \"None\" paragraph (the answer under the \"Prerequisite tasks\" heading)\n# had mistakenly been given the Heading 1 style with an explicit 11pt\n# (half-points value 22) size override, and carried a bookmark. Put it back\n# to plain Normal body text and drop the now-unneeded bookmark.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"None\" -and $p.Style.NameLocal -eq \"Heading 1\") {\n        $p.Range.Style = \"Normal\"\n        $p.Range.Font.Size = $null\n        break\n    }\n}\n\nif ($d.Bookmarks.Exists(\"_aqdfz55armzg\")) {\n    $d.Bookmarks.Item(\"_aqdfz55armzg\").Delete()\n}\n"}
